$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '46.889.86'
$ws.Range("E2").Value = '  +4.98%  '
$ws.Range("D3").Value = '2.354.00'
$ws.Range("E3").Value = '  +4.59%  '
$ws.Range("E4").Value = '  -0.70%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '307.88'
$ws.Range("E5").Value = '  +0.39%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '98.98'
$ws.Range("E6").Value = '  +4.10%  '
$ws.Range("E7").Value = '  +1.56%  '
$ws.Range("E8").Value = '  -0.62%  '
$ws.Range("E9").Value = '  +4.00%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '36.15'
$ws.Range("E10").Value = '  +3.52%  '
$ws.Range("E11").Value = '  +0.86%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.48'
$ws.Range("E12").Value = '  +3.55%  '
$ws.Range("E13").Value = '  -0.22%  '
$ws.Range("D14").Value = '2.708.21'
$ws.Range("E14").Value = '  +4.39%  '
$ws.Range("D15").Value = '2.342.08'
$ws.Range("E15").Value = '  +4.13%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.34'
$ws.Range("E16").Value = '  +5.36%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.835'
$ws.Range("E17").Value = '  +0.13%  '
$ws.Range("D18").Value = '46.747.89'
$ws.Range("E18").Value = '  +5.14%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.77'
$ws.Range("E19").Value = '  +17.11%  '
$ws.Range("D20").Value = '0.0₃0956'
$ws.Range("E20").Value = '  +1.88%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.22'
$ws.Range("E21").Value = '  +0.24%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '66.98'
$ws.Range("E22").Value = '  +2.44%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '245.65'
$ws.Range("E23").Value = '  +3.27%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.99'
$ws.Range("E24").Value = '  +1.01%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.00'
$ws.Range("E25").Value = '  +1.28%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("E26").Value = '  +0.00%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '42.08'
$ws.Range("E27").Value = '  +13.09%  '
$ws.Range("E28").Value = '  -0.13%  '
$ws.Range("E29").Value = '  +1.56%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '20.31'
$ws.Range("E30").Value = '  +1.57%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '5.80'
$ws.Range("E31").Value = '  -2.10%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '152.31'
$ws.Range("E32").Value = '  +2.31%  '
$ws.Range("E33").Value = '  +4.35%  '
$ws.Range("E34").Value = '  -0.20%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.16'
$ws.Range("E35").Value = '  -2.66%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.111'
$ws.Range("E36").Value = '  +2.32%  '
$ws.Range("E37").Value = '  +0.55%  '
$ws.Range("E38").Value = '  -0.83%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '4.11'
$ws.Range("E39").Value = '  +8.37%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0319'
$ws.Range("E40").Value = '  +6.28%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.47'
$ws.Range("E41").Value = '  +2.89%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '14.08'
$ws.Range("E42").Value = '  -7.98%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.999'
$ws.Range("E43").Value = '  -0.81%  '
$ws.Range("B44").Value = 'Maker'
$ws.Range("C44").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D44").Value = '1.860.65'
$ws.Range("E44").Value = '  +3.07%  '
$ws.Range("B45").Value = 'Stacks'
$ws.Range("C45").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.97'
$ws.Range("E45").Value = '  +11.16%  '
$ws.Range("E46").Value = '  +6.16%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '81.29'
$ws.Range("E47").Value = '  -0.80%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '73.76'
$ws.Range("E48").Value = '  +7.35%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '4.96'
$ws.Range("E49").Value = '  +2.41%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '99.72'
$ws.Range("E50").Value = '  +1.01%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '55.74'
$ws.Range("E51").Value = '  +3.01%  '
